# Function Tracker.xlsx - "All Bitbang Created, created generic comm function"
#
# 1. Insert a new row (new row 25) documenting the new generic SWD comm
#    function, pushing the existing Bitbang rows down by one.
# 2. Mark several functions as "Written" (column E) now that the
#    corresponding bitbang routines have been implemented.
# 3. Update the worksheet selection to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row for SWD_Comm -----------------------------------
$ws.Rows.Item(25).Insert()

$ws.Range("A25").Value = "uint32_t SWD_Comm(uint8_t command, uint32_t data)"
$ws.Range("B25").Value = "ProgrammerInterface.c"
$ws.Range("C25").Value = "Generic function to perform an SWD message"
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Formula = "=SUM(D25:G25)"

# --- Mark newly-written functions as Written (column E) ----------------
# setupSWDPins, configSWDPinsInput, configSWDPinsOutput, Clear_Target,
# SWD_Start
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("E22").Value = 1

# SWD_bitOut, SWD_bitIn, SWD_bitTurn (now rows 26-28 after the insert above)
$ws.Range("E26").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("E28").Value = 1

# --- Update view / selection state --------------------------------------
$ws.Range("L20").Select()

Write-Output "done"
